# Edit script for SC_unemployment.xlsx
# Applies: updated D-column values for rows 217-276 (revised JOLTS/LAUS data),
# new comments on D217:D276 ("Data were subject to revision..."), new rows
# 277-278 (Jan/Feb 2022 data) with the "Preliminary." comment moved to the
# new last row (D278), the "Years:" label bumped to "2000 to 2022", and the
# footer timestamp updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update revised values in column D for rows 217-276 ---------------
$revisedValues = @(
    @(217, 98609),
    @(218, 96842),
    @(219, 95309),
    @(220, 94216),
    @(221, 93629),
    @(222, 93696),
    @(223, 94370),
    @(224, 95210),
    @(225, 95747),
    @(226, 95652),
    @(227, 94504),
    @(228, 91832),
    @(229, 87840),
    @(230, 83236),
    @(231, 78871),
    @(232, 75584),
    @(233, 73997),
    @(234, 74031),
    @(235, 74961),
    @(236, 75978),
    @(237, 76399),
    @(238, 76058),
    @(239, 75386),
    @(240, 74935),
    @(241, 74755),
    @(242, 74427),
    @(243, 73554),
    @(244, 71570),
    @(245, 68313),
    @(246, 64378),
    @(247, 60681),
    @(248, 57974),
    @(249, 56620),
    @(250, 56876),
    @(251, 58438),
    @(252, 60933),
    @(253, 63929),
    @(254, 67120),
    @(255, 70081),
    @(256, 268537),
    @(257, 212235),
    @(258, 181338),
    @(259, 166767),
    @(260, 147497),
    @(261, 137252),
    @(262, 125520),
    @(263, 119089),
    @(264, 115021),
    @(265, 108934),
    @(266, 104414),
    @(267, 101358),
    @(268, 99319),
    @(269, 97924),
    @(270, 97345),
    @(271, 95396),
    @(272, 92837),
    @(273, 89382),
    @(274, 86523),
    @(275, 85167),
    @(276, 84737)
)

foreach ($pair in $revisedValues) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 4).Value() = $v
}

# --- 2. Add the revision comment to each of those same D217:D276 cells ---
# (AddComment replaces any prior comment on the cell, which is what we want
# for D276 -- it previously held the "Preliminary." note that now belongs
# on the new final data row, D278, added below.)
$revisionNote = "*  Data were subject to revision on March 2, 2022.`n"
for ($r = 217; $r -le 276; $r++) {
    $ws.Cells.Item($r, 4).AddComment($revisionNote)
}

# --- 3. Append the two new monthly rows (Jan 2022, Feb 2022) -------------
$ws.Range("A276:D276").Copy()
$ws.Range("A277:D278").PasteSpecial(-4122)

$ws.Cells.Item(277, 1).Value() = "LASST450000000000004"
$ws.Cells.Item(277, 2).Value() = 2022
$ws.Cells.Item(277, 3).Value() = "M01"
$ws.Cells.Item(277, 4).Value() = 82242

$ws.Cells.Item(278, 1).Value() = "LASST450000000000004"
$ws.Cells.Item(278, 2).Value() = 2022
$ws.Cells.Item(278, 3).Value() = "M02"
$ws.Cells.Item(278, 4).Value() = 82658

# --- 4. Move the "Preliminary." note onto the new last row, D278 ---------
$ws.Cells.Item(278, 4).AddComment("*  Preliminary.`n")

# --- 5. Bump the "Years:" summary label from 2000-2021 to 2000-2022 ------
$ws.Range("B10").Value() = "2000 to 2022"

# --- 6. Update the footer generation timestamp ----------------------------
$ws.PageSetup.LeftFooter = "Source: Bureau of Labor Statistics"
$ws.PageSetup.RightFooter = "Generated on: March 28, 2022 (06:22:58 PM)"
